$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A41").Value = 0
$ws.Range("B41").Value = 0
$ws.Range("C41").Value = 0.256281
$ws.Range("D41").Value = -0.1729463799634047
$ws.Range("E41").Value = "query"
